$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.110.82"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.750.63"
$ws.Range("E3").Value = "  +3.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.61"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.04"
$ws.Range("E6").Value = "  +4.98%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.547"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.750.77"
$ws.Range("E9").Value = "  +3.78%  "
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.364"
$ws.Range("E11").Value = "  +3.43%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.35"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.157"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.89"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.251.23"
$ws.Range("E15").Value = "  +3.80%  "
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.037.15"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.733.45"
$ws.Range("E18").Value = "  +3.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.00"
$ws.Range("E19").Value = "  +5.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.76"
$ws.Range("E20").Value = "  +5.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "369.22"
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("E23").Value = "  +3.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.05"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.99"
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.877.12"
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "601.56"
$ws.Range("E30").Value = "  +7.47%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("E32").Value = "  +4.01%  "
$ws.Range("E33").Value = "  +3.93%  "
$ws.Range("E34").Value = "  +5.95%  "
$ws.Range("E35").Value = "  +3.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.64"
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.40"
$ws.Range("E38").Value = "  +2.44%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.20"
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.383"
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("E42").Value = "  +3.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.73"
$ws.Range("E43").Value = "  +3.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.04"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("E45").Value = "  -4.47%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "159.32"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.98"
$ws.Range("E48").Value = "  +5.59%  "
$ws.Range("E49").Value = "  +6.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.612"
$ws.Range("E50").Value = "  +7.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.20"
$ws.Range("E51").Value = "  -0.27%  "
